$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the CL1 scores for students (column C, rows 15-23)
$ws.Range("C15").Value = 9
$ws.Range("C16").Value = 10
$ws.Range("C17").Value = 10
$ws.Range("C18").Value = 10
$ws.Range("C19").Value = 10
$ws.Range("C20").Value = 9
$ws.Range("C21").Value = 9
$ws.Range("C22").Value = 10
$ws.Range("C23").Value = 9

# Row 24 gets explicit zeros across CL1-CL4
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0

# Update the view: scroll so row 4 is the top-left row and select C15
$ws.Range("C15").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
